$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data range (rows 2-31) before writing the refreshed dataset.
$ws.Range("A2:F31").ClearContents()

# Column B holds dates formatted as literal text ("01/01/2015"); format the
# target range as Text first so Excel doesn't auto-convert them to date serials.
$ws.Range("B2:B34").NumberFormat = "@"

# New dataset now spans rows 2-34 (33 data rows) after the manual figure refresh
# (adds 01/01/2025 rows for Brasil/Nordeste, reorders/extends the Sergipe block).
$data = New-Object "object[,]" 33,6

$data[0,0]="Brasil"; $data[0,1]="01/01/2015"; $data[0,2]="Feminicídio"; $data[0,3]=0.6273661320339025; $data[0,4]=$null; $data[0,5]=$true
$data[1,0]="Brasil"; $data[1,1]="01/01/2016"; $data[1,2]="Feminicídio"; $data[1,3]=0.863032297756537; $data[1,4]=$null; $data[1,5]=$true
$data[2,0]="Brasil"; $data[2,1]="01/01/2017"; $data[2,2]="Feminicídio"; $data[2,3]=1.082538168435736; $data[2,4]=$null; $data[2,5]=$true
$data[3,0]="Brasil"; $data[3,1]="01/01/2018"; $data[3,2]="Feminicídio"; $data[3,3]=1.236108881648152; $data[3,4]=$null; $data[3,5]=$true
$data[4,0]="Brasil"; $data[4,1]="01/01/2019"; $data[4,2]="Feminicídio"; $data[4,3]=1.49246755203831; $data[4,4]=$null; $data[4,5]=$false
$data[5,0]="Brasil"; $data[5,1]="01/01/2020"; $data[5,2]="Feminicídio"; $data[5,3]=1.519724099657912; $data[5,4]=$null; $data[5,5]=$false
$data[6,0]="Brasil"; $data[6,1]="01/01/2021"; $data[6,2]="Feminicídio"; $data[6,3]=1.570641743846045; $data[6,4]=$null; $data[6,5]=$false
$data[7,0]="Brasil"; $data[7,1]="01/01/2022"; $data[7,2]="Feminicídio"; $data[7,3]=1.538407276001129; $data[7,4]=$null; $data[7,5]=$false
$data[8,0]="Brasil"; $data[8,1]="01/01/2023"; $data[8,2]="Feminicídio"; $data[8,3]=1.511361205446982; $data[8,4]=$null; $data[8,5]=$false
$data[9,0]="Brasil"; $data[9,1]="01/01/2024"; $data[9,2]="Feminicídio"; $data[9,3]=1.450683841983206; $data[9,4]=$null; $data[9,5]=$false
$data[10,0]="Brasil"; $data[10,1]="01/01/2025"; $data[10,2]="Feminicídio"; $data[10,3]=0.86356972985822; $data[10,4]=$null; $data[10,5]=$false
$data[11,0]="Nordeste"; $data[11,1]="01/01/2015"; $data[11,2]="Feminicídio"; $data[11,3]=0.6635307538337432; $data[11,4]=$null; $data[11,5]=$true
$data[12,0]="Nordeste"; $data[12,1]="01/01/2016"; $data[12,2]="Feminicídio"; $data[12,3]=0.9129393838946243; $data[12,4]=$null; $data[12,5]=$true
$data[13,0]="Nordeste"; $data[13,1]="01/01/2017"; $data[13,2]="Feminicídio"; $data[13,3]=1.239969626186528; $data[13,4]=$null; $data[13,5]=$true
$data[14,0]="Nordeste"; $data[14,1]="01/01/2018"; $data[14,2]="Feminicídio"; $data[14,3]=1.326523952625601; $data[14,4]=$null; $data[14,5]=$true
$data[15,0]="Nordeste"; $data[15,1]="01/01/2019"; $data[15,2]="Feminicídio"; $data[15,3]=1.497286779739304; $data[15,4]=$null; $data[15,5]=$false
$data[16,0]="Nordeste"; $data[16,1]="01/01/2020"; $data[16,2]="Feminicídio"; $data[16,3]=1.421470954921448; $data[16,4]=$null; $data[16,5]=$false
$data[17,0]="Nordeste"; $data[17,1]="01/01/2021"; $data[17,2]="Feminicídio"; $data[17,3]=1.456790721661446; $data[17,4]=$null; $data[17,5]=$false
$data[18,0]="Nordeste"; $data[18,1]="01/01/2022"; $data[18,2]="Feminicídio"; $data[18,3]=1.338905345285876; $data[18,4]=$null; $data[18,5]=$false
$data[19,0]="Nordeste"; $data[19,1]="01/01/2023"; $data[19,2]="Feminicídio"; $data[19,3]=1.327731384986638; $data[19,4]=$null; $data[19,5]=$false
$data[20,0]="Nordeste"; $data[20,1]="01/01/2024"; $data[20,2]="Feminicídio"; $data[20,3]=1.353249472239435; $data[20,4]=$null; $data[20,5]=$false
$data[21,0]="Nordeste"; $data[21,1]="01/01/2025"; $data[21,2]="Feminicídio"; $data[21,3]=0.7961563668809312; $data[21,4]=$null; $data[21,5]=$false
$data[22,0]="Sergipe"; $data[22,1]="01/01/2015"; $data[22,2]="Feminicídio"; $data[22,3]=0; $data[22,4]=20.5; $data[22,5]=$true
$data[23,0]="Sergipe"; $data[23,1]="01/01/2016"; $data[23,2]="Feminicídio"; $data[23,3]=0; $data[23,4]=23; $data[23,5]=$true
$data[24,0]="Sergipe"; $data[24,1]="01/01/2017"; $data[24,2]="Feminicídio"; $data[24,3]=1.715876576997817; $data[24,4]=5; $data[24,5]=$true
$data[25,0]="Sergipe"; $data[25,1]="01/01/2018"; $data[25,2]="Feminicídio"; $data[25,3]=1.359138849624878; $data[25,4]=13; $data[25,5]=$true
$data[26,0]="Sergipe"; $data[26,1]="01/01/2019"; $data[26,2]="Feminicídio"; $data[26,3]=1.766753237575308; $data[26,4]=7; $data[26,5]=$false
$data[27,0]="Sergipe"; $data[27,1]="01/01/2020"; $data[27,2]="Feminicídio"; $data[27,3]=1.166870869068754; $data[27,4]=20; $data[27,5]=$false
$data[28,0]="Sergipe"; $data[28,1]="01/01/2021"; $data[28,2]="Feminicídio"; $data[28,3]=1.651913700724447; $data[28,4]=9; $data[28,5]=$false
$data[29,0]="Sergipe"; $data[29,1]="01/01/2022"; $data[29,2]="Feminicídio"; $data[29,3]=1.555611048113412; $data[29,4]=12; $data[29,5]=$false
$data[30,0]="Sergipe"; $data[30,1]="01/01/2023"; $data[30,2]="Feminicídio"; $data[30,3]=1.298940632730219; $data[30,4]=19; $data[30,5]=$false
$data[31,0]="Sergipe"; $data[31,1]="01/01/2024"; $data[31,2]="Feminicídio"; $data[31,3]=0.8052392083854389; $data[31,4]=26; $data[31,5]=$false
$data[32,0]="Sergipe"; $data[32,1]="01/01/2025"; $data[32,2]="Feminicídio"; $data[32,3]=0.5592613276384552; $data[32,4]=24; $data[32,5]=$false

$ws.Range("A2:F34").Value2 = $data
